$d = $word.ActiveDocument

# 1. "Mom (neutral surprised): How was it?" -> "Mom (neutral smiling): How was it?"
$d.Content.Find.Execute("Mom (neutral surprised): How was it?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mom (neutral smiling): How was it?", 2)

# 2. "Mom (neutral smiling): I see…" -> "Mom (neutral smiling_nervous): I see…"
$d.Content.Find.Execute("Mom (neutral smiling): I see", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mom (neutral smiling_nervous): I see", 2)

# 3. Insert a new paragraph "Mom (exit):" right after the "Pro: Oh, thanks." paragraph.
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Pro: Oh, thanks.") {
        $targetIdx = $idx
        $p.Range.InsertParagraphAfter()
        break
    }
}
if ($targetIdx -gt 0) {
    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.Text = "Mom (exit):"
}
